$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) slide49 ("Background: Usenet"): drop the trailing empty paragraph that
#    used to sit after "Articles shared via flood-fill" in the content
#    placeholder.
# ---------------------------------------------------------------------------
$bgSlide = $p.Slides.Item(49)
$bgContent = $bgSlide.Shapes.Item("Content Placeholder 3")
$bgRange = $bgContent.TextFrame.TextRange
$trailingPara = $bgRange.Paragraphs(9, 1)
if ($trailingPara.Text -eq "") {
    $trailingPara.Delete()
}

# ---------------------------------------------------------------------------
# 2) Add the new "UsenetDHT" slide right after it (slide 50), by duplicating
#    slide49 (keeps the title/sldNum/content placeholder layout, ids and
#    timing/colour-map overrides identical to its sibling slides) and then
#    replacing its shapes/text.
# ---------------------------------------------------------------------------
$bgSlide.Duplicate() | Out-Null
$newSlide = $p.Slides.Item(50)

# Drop the picture and caption textbox that were copied from slide49 - the
# new slide only has Title / Slide Number / Content placeholders.
$newSlide.Shapes.Item("Picture 5").Delete()
$newSlide.Shapes.Item("TextBox 6").Delete()

# --- Title -------------------------------------------------------------
$title = $newSlide.Shapes.Item("Title 1")
$title.TextFrame.TextRange.Text = "UsenetDHT"

# --- Content placeholder -------------------------------------------------
$content = $newSlide.Shapes.Item("Content Placeholder 3")
$tr = $content.TextFrame.TextRange

$tr.Text = "Problem:`rEach server stores copies of all articles (that it wants)`rO(n) copies of each article!`rIdea:`rStore articles in common store`rO(n) reduction of space used`rUsenetDHT:`rPeer-to-peer applications`rEach node acts as Usenet frontend, and DHT node`rHeaders flood-filled as normal, articles stored in DHT"

# Paragraph indent levels (COM IndentLevel is 1-based; lvl=1 -> no/0 pPr).
$tr.Paragraphs(2, 1).IndentLevel = 2
$tr.Paragraphs(3, 1).IndentLevel = 2
$tr.Paragraphs(4, 1).IndentLevel = 1
$tr.Paragraphs(5, 1).IndentLevel = 2
$tr.Paragraphs(6, 1).IndentLevel = 2
$tr.Paragraphs(7, 1).IndentLevel = 1
$tr.Paragraphs(8, 1).IndentLevel = 2
$tr.Paragraphs(9, 1).IndentLevel = 2
$tr.Paragraphs(10, 1).IndentLevel = 2

# Split the two "two-run" paragraphs so the trailing clause gets its own
# run (mirrors the source deck's baseline="0" runs from autocorrect).
$run2 = $tr.Paragraphs(2, 1).Characters(29, 30)
$run2.Font.BaselineOffset = 0

$run9 = $tr.Paragraphs(9, 1).Characters(25, 23)
$run9.Font.BaselineOffset = 0
